# Fix for date/time values that include seconds: the fractional-day value
# now gets written out with full (rounded) precision, e.g. 11/5/2013 11:45:00
# instead of rounding down to 11:44:59. Demonstrate the fix with a new test
# row holding the previously-mis-rounded timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 41583.489583333299 == 11/5/2013 11:45:00 (stored as a fractional day count);
# this is the exact value that used to round down to 11:44:59 before the fix.
$dateValue = 41583.489583333299

# Build the new custom style (number format + wrapped text) on A8 first ...
$a8 = $ws.Range("A8")
$a8.Value = $dateValue
$a8.NumberFormat = "m/d/yyyy\ h:mm:ss;@"
$a8.WrapText = $true

# ... then stamp the same value onto B8/C8 and clone A8's format onto them,
# so all three cells end up sharing one single cell style.
$b8 = $ws.Range("B8")
$b8.Value = $dateValue
$c8 = $ws.Range("C8")
$c8.Value = $dateValue

$a8.Copy() | Out-Null
$ws.Range("B8:C8").PasteSpecial(-4122) | Out-Null

# Move the active selection (matches the workbook's last saved cursor spot).
$ws.Range("C12").Select() | Out-Null
